$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room for more time-log rows: insert 4 blank rows right
#    before the old row 67 (the row holding the "TODO-- PW feedback
#    authentication..." note). This pushes everything from the old
#    row 67 downward by 4 rows while preserving all of its content
#    and formatting.
# ------------------------------------------------------------------
$ws.Rows("67:70").Insert()

# ------------------------------------------------------------------
# 2) Rows 64 and 65 held two older notes ("Note: I may have needed
#    to do something..." and the tutorials.jenkov.com link). They
#    sit above the insertion point, so the insert above did not move
#    them. In the final layout they live a little further down (rows
#    68 and 69), so relocate them there now, before those two rows
#    are reused for new content.
# ------------------------------------------------------------------
$ws.Cells.Item(64, 4).Cut($ws.Cells.Item(68, 4))
$ws.Cells.Item(65, 4).Cut($ws.Cells.Item(69, 4))
$ws.Cells.Item(65, 4).Clear()

# ------------------------------------------------------------------
# 3) Fill in the newly-recovered time log entries.
# ------------------------------------------------------------------
$ws.Cells.Item(57, 1).Value = 43562
$ws.Cells.Item(57, 2).Value = 5
$ws.Cells.Item(57, 4).Value = "Team Project: looked at what Kelly did (baby web app)`nWeek 10: created branch and pull request`nWeek 9: worked on activity"
$ws.Rows(57).RowHeight = 45

$ws.Cells.Item(58, 1).Value = 43563
$ws.Cells.Item(58, 2).Value = 1
$ws.Cells.Item(58, 4).Value = "Finished Week 9 activity"

$ws.Cells.Item(59, 1).Value = 43564
$ws.Cells.Item(59, 2).Value = 4.5
$ws.Cells.Item(59, 4).Value = "Team project: tried to get log4j out of github; added a service method that provides a party parrot based on name; tried to research some way of handling json info more directly, but did not figure something out so ended up working with a list of all the parrots, as objects, for my method.  `nIssue:Time estimate may be wrong - I lost some time log data while going back and forth between branch and  master of my indie project."
$ws.Rows(59).RowHeight = 75

$ws.Cells.Item(60, 1).Value = 43566

$ws.Cells.Item(62, 4).Value = "early AM - 1 hour, + 8:35 - x"
$ws.Cells.Item(63, 4).Value = "created properties file - need to make application upload and use it"
$ws.Cells.Item(64, 4).Value = "populated jsp with ""top"" careers to search"

# ------------------------------------------------------------------
# 4) Leave the selection where the author ended up working.
# ------------------------------------------------------------------
$ws.Range("D65").Select()
